# Updated symbol list on Fri Dec 30 03:24:19 UTC 2022 with GitHub Actions
# Refreshes coin prices/volume labels in the "cryptos" sheet, and re-syncs
# several rows whose coin ranking shifted position.
# Note: price values in column D are stored as text (matching the source
# data, which preserves formatting such as trailing zeros), so a leading
# apostrophe is used to force literal text instead of auto-converted numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''245.72'

$ws.Range("D3").Value = '''24.22'

$ws.Range("D5").Value = '''0.05784'

$ws.Range("D6").Value = '''6.495'

$ws.Range("D7").Value = '''3.153'

$ws.Range("D8").Value = '''0.8181'

$ws.Range("D9").Value = '''0.8501'

$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '''0.009813'
$ws.Range("E10").Value = '9OneONEBestin24h'

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1360'
$ws.Range("E11").Value = '10WazirXWRX'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.06953'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'

$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = '''0.03155'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("D14").Value = '''0.02880'

$ws.Range("D15").Value = '''0.09385'

$ws.Range("D16").Value = '''3.749'

$ws.Range("D17").Value = '''0.001512'

$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D19").Value = '''0.006282'
$ws.Range("E19").Value = '18TigerCashTCH'

$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D20").Value = '''0.001239'
$ws.Range("E20").Value = '19BitKanKAN'

$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = '''0.004607'
$ws.Range("E21").Value = '20HotbitTokenHTB'

$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").Value = '''0.00006902'
$ws.Range("E22").Value = '21NitroExNTXWorstin24h'

$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = '''3.521'
$ws.Range("E23").Value = '22LEOLEO'

$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = '''2.148'
$ws.Range("E24").Value = '23BTSETokenBTSE'

$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").Value = '''0.3191'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'

$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").Value = '''0.1348'
$ws.Range("E26").Value = '25ProBitTokenPROB'

$ws.Range("D27").Value = '''0.1327'

$ws.Range("D28").Value = '''0.0002330'

$ws.Range("D40").Value = '''0.03653'

$ws.Range("D41").Value = '''0.006242'

$ws.Range("D42").Value = '''0.1054'

$ws.Range("D43").Value = '''0.003401'

$ws.Range("D44").Value = '''0.007451'

$ws.Range("D45").Value = '''0.00005277'

$ws.Range("D47").Value = '''0.3612'

$ws.Range("D48").Value = '''0.002332'
$ws.Range("E48").Value = '47BOLOBOLO'

$ws.Range("D49").Value = '''0.00002101'

$ws.Range("D50").Value = '''0.0002001'
